$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.250.57'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.72%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.855.41'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.19%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9998'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.08%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.18'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.93%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9997'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.11%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4797'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.10%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2816'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.13%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06502'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.13%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.856.35'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.10%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07339'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.29%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '16.30'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.20%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.136'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.10%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '87.41'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.40%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6488'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.68%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '30.211.10'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.73%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.31'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.19%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9996'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.03%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007655'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.42%  '

# Row 20
$ws.Range("B20").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C20").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.102.26'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.74%  '

# Row 21
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.304'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.86%  '

# Row 22
$ws.Range("B22").Value = 'BitcoinCash'
$ws.Range("C22").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '222.57'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +16.54%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.14%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.085'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.89%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.272'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.24%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.73'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.04%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.52'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.41%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.927'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.47%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.436'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.37%  '

# Row 30
$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09206'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.06%  '

# Row 31
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.245'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.85%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.974'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.75%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05024'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.17%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7418'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.11%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.151'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.37%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.686'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.95%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01822'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.88%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.617'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.00%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.9082'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.52%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.061'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.57%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.979'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.68%  '

# Row 42
$ws.Range("E42").Value = '  +0.12%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4263'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.02%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9998'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.68%  '

# Row 45
$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.422'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.30%  '

# Row 46
$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1322'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.51%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.557'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +10.65%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '64.30'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.63%  '

# Row 49
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.782'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.38%  '

# Row 50
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.28'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.15%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05665'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.67%  '
